$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current table ends at row 8 (the "closing" row with thick bottom border).
# We need to insert a new row before it (so it becomes the new row 8),
# pushing the old closing row down to row 9, then put the new request
# text into the new row 8... but per the diff, actually the NEW content
# row ends up as the last row (row 9) and keeps the "closing" style,
# while what used to be the closing row (row 8, "Supprimer...") becomes
# a normal middle row.
#
# Simplest way to reproduce this: select the last row of the table
# (row 8, the thick-bottom "closing" row) and insert a new row ABOVE it
# by copying its format, then restore the old row 8 to a normal middle
# row style, and fill the new row with the new text.

# Insert a new row at row 9 position: insert before row 8, shifting old
# row 8 down to row 9. Excel's native "insert copied cells" / "insert row"
# behavior adjusts borders similarly to what's in the diff.

$ws.Rows.Item(8).Insert(-4121)  # xlShiftDown = -4121

# Copy the style/format from the (now) row 9 thick-bottom row into row 8,
# so row 8 reverts to being a "normal middle" row like rows 3-7, and the
# newly inserted blank row 8 gets the normal middle style too.
$ws.Range("A7:C7").Copy() | Out-Null
$ws.Range("A8:C8").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Now give the last row (row 9, which still has old row 8's original
# closing format because we inserted above it, not copied over it) the
# text content of the new request.
$ws.Range("A9").Value = "Il faudrait rajouter la possibilité de mettre du texte pour les Daily rate pour qu’on puisse mettre le type de monnaie"

$excel.CutCopyMode = 0

$ws.Range("A17").Select() | Out-Null
